# Design Fire Creator - Validation.xlsx update
# - Move/resize the "Diagram 1" chart graphic frame on Ark1
# - Clear two stray debug formulas (S11, R15) on Ark1
# - Change the active selection on Ark1 to W5 (and drop the old top-left/tab-selected state)
# - Add a new "Test" worksheet (UnitTester) after Ark1, with two calculation blocks
#   (Growth phase + Decay phase) and make it the active sheet

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Reposition the "Diagram 1" chart (Plotly-style growth/decay chart) ---
foreach ($co in $ws1.ChartObjects()) {
    if ($co.Name -eq "Diagram 1") {
        $co.Left = 1326.2284645669292
        $co.Top = 250.44826771653544
        $co.Width = 518.2032480314961
        $co.Height = 283.46456692913387
    }
}

# --- Clear the stray leftover debug formulas ---
$ws1.Range("S11").ClearContents()
$ws1.Range("R15").ClearContents()

# --- Update selection on Ark1 (also drops tabSelected / topLeftCell once Test becomes active) ---
$ws1.Range("W5").Select() | Out-Null

# --- Add the new "Test" sheet right after Ark1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Test"

# --- Growth Phase block ---
$ws2.Range("B3").Value = "Growth Phase with known Duration and Growth Rate Factor"
$ws2.Range("B3").Font.Bold = $true

$ws2.Range("B5").Value = "Duration "
$ws2.Range("C5").Value = 100
$ws2.Range("B6").Value = "Growth Rate Factor"
$ws2.Range("C6").Value = 0.047
$ws2.Range("D6").Value = "kW/s²"
$ws2.Range("D5").Value = "s"

$ws2.Range("B8").Value = "Effect"
$ws2.Range("C8").Formula = "=C6*C5^2"
$ws2.Range("D8").Value = "kW"

# --- Decay Phase block ---
$ws2.Range("B12").Value = "Decay Phase with known Duration and Growth Rate Factor"
$ws2.Range("B12").Font.Bold = $true

$ws2.Range("C14").Value = 100
$ws2.Range("B15").Value = "InitialYq"
$ws2.Range("C15").Value = 1000
$ws2.Range("B14").Value = "InitialXt"

$ws2.Range("B17").Value = "Duration "
$ws2.Range("C17").Value = 100
$ws2.Range("B18").Value = "Growth Rate Factor"
$ws2.Range("C18").Value = 0.047
$ws2.Range("D18").Value = "kW/s²"
$ws2.Range("D17").Value = "s"

$ws2.Range("B20").Value = "Effect"
$ws2.Range("C20").Formula = "=C15-C18*C17^2"

# Column B needs to be wide enough for "Growth Rate Factor" (the widest
# non-title label) -- matches the width produced by autofitting the column
# to that text.
$ws2.Columns.Item(2).ColumnWidth = 20.0

$ws2.Range("D26").Select() | Out-Null
